# Preparation publication 0.2.0
# - Insert a new "Jurisdiction | iso:code:3166:FR" row into the Metadata sheet
#   (right after the "Contact" row, before "Description")
# - Bump Version 0.1.1 -> 0.2.0
# - Bump Date 2023-10-20T07:19:33+00:00 -> 2023-10-20T08:59:58+00:00

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row 11 (pushing Description/Purpose/Copyright/Source/Target down by one)
$ws.Rows.Item(11).Insert()

# Copy formatting (fill/border/font/alignment) from the row above (Contact row)
# so the new row matches the rest of the table's style instead of Excel's
# default insert formatting.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new Jurisdiction row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Update Version and Date values
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

Write-Output "done"
